# Scenario 3a/3b coverage parameter adjustments (Platform Coverage sheet)
#
# Summary of intent (per commit message "Adjustments to scenario 3a and 3b
# parameters for STH"):
#   - The MDA/Treatment/Campaign age band 5-15 (row 2) keeps its early-year
#     coverage (2018/2020/2022/2024 = 0.6) but loses its later-year
#     (2026-2040) coverage values.
#   - A brand-new MDA/Treatment/Campaign row is inserted for age band 2-15
#     with 0.8 coverage for 2026-2040 (even years only).
#   - The existing 15-50 and 50-65 age bands gain 0.5 coverage for
#     2026-2040 (even years only).
#   - The Vaccine rows (EPI / School / Out-of-school campaign) are
#     unchanged in content, just shifted down one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Platform Coverage")

$years = @("P","R","T","V","X","Z","AB","AD")

# --- Row 2 (MDA, Treatment, Campaign, age 5-15): drop the 2026-2040 values,
#     keep only the 2018-2024 ones already present. ---
foreach ($col in $years) {
    $ws.Range($col + "2").ClearContents()
}

# --- Insert a brand-new row 3 for the MDA / Treatment / Campaign age
#     band 2-15, pushing the old rows 3-7 down to 4-8. ---
$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value = "All"
$ws.Range("B3").Value = "Treatment"
$ws.Range("C3").Value = "Campaign"
$ws.Range("D3").Value = "MDA"
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 15
foreach ($col in $years) {
    $ws.Range($col + "3").Value = 0.8
}

# --- Row 4 (formerly row 3: MDA, Treatment, Campaign, age 15-50): add the
#     new 2026-2040 coverage values (0.5). ---
foreach ($col in $years) {
    $ws.Range($col + "4").Value = 0.5
}

# --- Row 5 (formerly row 4: MDA, Treatment, Campaign, age 50-65): add the
#     new 2026-2040 coverage values (0.5). ---
foreach ($col in $years) {
    $ws.Range($col + "5").Value = 0.5
}

# Rows 6-8 (formerly rows 5-7: Vaccine / EPI / School / Out-of-school
# campaign) are carried down unchanged by the row insert above - no
# further edits needed for them.

# --- Sheet view cosmetics: zoom level and active selection. ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 120
$ws.Range("AD2").Select()
